$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("cor_z")
$ws.Range("B3").Value = -32.16089653357908
$ws.Range("B4").Value = -29.0368165335791
$ws.Range("B5").Value = -25.91273653357869
$ws.Range("B6").Value = -22.78865653357826
$ws.Range("B7").Value = -19.66457653357855
$ws.Range("B8").Value = -16.540496533578
$ws.Range("B9").Value = -13.41641653357864
$ws.Range("B10").Value = -10.2923365335791
$ws.Range("B11").Value = -7.168256533578247
$ws.Range("B12").Value = -4.044176533578007
$ws.Range("B13").Value = -0.9200965335777663
$ws.Range("B14").Value = 2.203983466421769
$ws.Range("B15").Value = 5.328063466421129
$ws.Range("B16").Value = 8.452143466421722
$ws.Range("B17").Value = 11.57622346642161
$ws.Range("B18").Value = 14.70030346642097
$ws.Range("B19").Value = 17.82438346642156
$ws.Range("B20").Value = 20.94846346642127
$ws.Range("B21").Value = 24.07254346642043

$ws = $wb.Worksheets.Item("mome_y")
$ws.Range("B3").Value = 4117.574222658763
$ws.Range("B4").Value = 2922.65827616597
$ws.Range("B5").Value = 1849.740777753168
$ws.Range("B6").Value = 898.8217274204068
$ws.Range("B7").Value = 69.90112516763834
$ws.Range("B8").Value = -637.0210290051448
$ws.Range("B9").Value = -1221.944735097871
$ws.Range("B10").Value = -1684.869993110704
$ws.Range("B11").Value = -2025.796803043466
$ws.Range("B12").Value = -2244.725164896232
$ws.Range("B13").Value = -2341.655078668975
$ws.Range("B14").Value = -2316.586544361723
$ws.Range("B15").Value = -2169.519561974502
$ws.Range("B16").Value = -1900.4541315073
$ws.Range("B17").Value = -1509.390252960033
$ws.Range("B18").Value = -996.327926332839
$ws.Range("B19").Value = -361.2671516256222
$ws.Range("B20").Value = 395.7920711616177
$ws.Range("B21").Value = 1274.849742028812

$ws = $wb.Worksheets.Item("drz")
$ws.Range("B2").Value = -0.00000000000000001219024528012889
$ws.Range("B3").Value = 0.08774797512503844
$ws.Range("B4").Value = 0.3084647816943936
$ws.Range("B5").Value = 0.6235629969956872
$ws.Range("B6").Value = 0.9983948944099221
$ws.Range("B7").Value = 1.402252443411484
$ws.Range("B8").Value = 1.808367309568142
$ws.Range("B9").Value = 2.193910854541042
$ws.Range("B10").Value = 2.539994136084718
$ws.Range("B11").Value = 2.83166790804708
$ws.Range("B12").Value = 3.057922620369422
$ws.Range("B13").Value = 3.21168841908642
$ws.Range("B14").Value = 3.289835146326132
$ws.Range("B15").Value = 3.293172340309997
$ws.Range("B16").Value = 3.226449235352837
$ws.Range("B17").Value = 3.098354761862852
$ws.Range("B18").Value = 2.921517546341628
$ws.Range("B19").Value = 2.712505911384128
$ws.Range("B20").Value = 2.491827875678701
$ws.Range("B21").Value = 2.283931154007075
